$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 7; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $colLetter = [char](64 + $c)
        $cellRef = "$colLetter$r"
        $ws.Range($cellRef).Value = $cellRef
    }
}

$ws.Range("C8").Select()
